$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.451.70'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''315.15'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '''0.5144'
$ws.Range("E7").Value = '  -3.24%  '
$ws.Range("D8").Value = '''0.3915'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = '''0.07693'
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("D10").Value = '''41.83'
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("E12").Value = '  +3.05%  '
$ws.Range("D13").Value = '''6.285'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").Value = '''7.549'
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").Value = '1.826.40'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '''93.61'
$ws.Range("E17").Value = '  +5.03%  '
$ws.Range("D18").Value = '''0.00001108'
$ws.Range("E18").Value = '  +3.85%  '
$ws.Range("D19").Value = '''0.06707'
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").Value = '''17.71'
$ws.Range("E20").Value = '  +2.30%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").Value = '''6.195'
$ws.Range("E22").Value = '  +2.98%  '
$ws.Range("D23").Value = '28.479.85'
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").Value = '''11.16'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").Value = '''2.255'
$ws.Range("E25").Value = '  +7.96%  '
$ws.Range("D26").Value = '''156.94'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").Value = '''20.63'
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("D28").Value = '2.038.62'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '''124.63'
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("D31").Value = '''1.119'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '''5.675'
$ws.Range("E33").Value = '  +2.29%  '
$ws.Range("D34").Value = '''3.657'
$ws.Range("E34").Value = '  -0.28%  '
$ws.Range("D35").Value = '''0.07010'
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("D36").Value = '''0.2219'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '''8.952'
$ws.Range("E37").Value = '  +5.27%  '
$ws.Range("D38").Value = '''0.02326'
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("D39").Value = '''5.152'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("D41").Value = '''11.23'
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").Value = '''1.181'
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").Value = '''1.393'
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").Value = '''0.5903'
$ws.Range("E46").Value = '  +2.73%  '
$ws.Range("D47").Value = '''3.710'
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("D48").Value = '''124.91'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '''1.979'
$ws.Range("E49").Value = '  +2.31%  '
$ws.Range("D50").Value = '''1.200'
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").Value = '''0.06926'
$ws.Range("E51").Value = '  +1.67%  '
